$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "exp": insert three new experiment rows (lowO2 / midO2 / highO2) just
# before the moshammer row, keeping the outcome/plot/plot columns as-is.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("exp")

# Row 8 currently holds "moshammer_2016_dme_mod.xlsx"; push it down by
# inserting 3 blank rows above it.
$ws1.Rows.Item(8).Resize(3).Insert() | Out-Null

$ws1.Range("A8").Value = "couch_2022_dme_lowO2.xlsx"
$ws1.Range("B8").Value = "outcome"
$ws1.Range("C8").Value = "plot"
$ws1.Range("D8").Value = "plot"

$ws1.Range("A9").Value = "couch_2022_dme_midO2.xlsx"
$ws1.Range("B9").Value = "outcome"
$ws1.Range("C9").Value = "plot"
$ws1.Range("D9").Value = "plot"

$ws1.Range("A10").Value = "couch_2022_dme_highO2.xlsx"
$ws1.Range("B10").Value = "outcome"
$ws1.Range("C10").Value = "plot"
$ws1.Range("D10").Value = "plot"

# ---------------------------------------------------------------------------
# Sheet "mech": collapse the three candidate mechanisms down to a single
# baseline mechanism.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("mech")

# Remove the v1 (Burke+Cl) and v3 (w/Couch mods) rows, leaving only the v2
# row, which becomes the new "Baseline" mechanism.
$ws2.Rows.Item(4).Delete() | Out-Null
$ws2.Rows.Item(3).Delete() | Out-Null

$ws2.Range("A2").Value = "dme_couch_v2.cti"
$ws2.Range("B2").Value = "dme_couch.csv"
$ws2.Range("C2").Value = "Baseline"

# ---------------------------------------------------------------------------
# View-state touch-ups captured in the workbook (selection + show formulas).
# Apply the "mech" selection first, then finish on "exp" so it ends up the
# active/selected sheet, matching the saved file.
# ---------------------------------------------------------------------------
$ws2.Range("A2:C2").Select() | Out-Null

$ws1.Range("C15").Select() | Out-Null
$excel.ActiveWindow.DisplayFormulas = $true
